# Auto-applies the cryptos.xlsx price/volume/coin-order refresh described
# in the commit "Updated cryptos list ... with GitHub Actions".
# A leading apostrophe forces each assignment to be stored as literal text
# (matching the inlineStr cells in the workbook) instead of letting Excel
# auto-coerce number-like strings (e.g. "1.000", "7.067") into numerics,
# which would silently drop meaningful trailing/format digits.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.235.78"
$ws.Range("E2").Value = "'  -0.67%  "

# Row 3
$ws.Range("D3").Value = "'1.861.60"
$ws.Range("E3").Value = "'  -0.84%  "

# Row 4
$ws.Range("E4").Value = "'  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'0.7147"
$ws.Range("E5").Value = "'  -0.85%  "

# Row 6
$ws.Range("D6").Value = "'240.79"
$ws.Range("E6").Value = "'  +0.23%  "

# Row 7
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "'  +0.07%  "

# Row 8
$ws.Range("B8").Value = "'Dogecoin"
$ws.Range("C8").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "'0.07722"
$ws.Range("E8").Value = "'  -1.50%  "

# Row 9
$ws.Range("B9").Value = "'Cardano"
$ws.Range("C9").Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.3083"
$ws.Range("E9").Value = "'  -0.80%  "

# Row 10
$ws.Range("D10").Value = "'24.90"
$ws.Range("E10").Value = "'  -0.41%  "

# Row 11
$ws.Range("D11").Value = "'0.08318"
$ws.Range("E11").Value = "'  +0.72%  "

# Row 12
$ws.Range("D12").Value = "'1.869.79"
$ws.Range("E12").Value = "'  +0.00%  "

# Row 13
$ws.Range("D13").Value = "'0.7171"
$ws.Range("E13").Value = "'  -1.48%  "

# Row 14
$ws.Range("D14").Value = "'5.219"
$ws.Range("E14").Value = "'  -1.22%  "

# Row 15
$ws.Range("D15").Value = "'90.87"
$ws.Range("E15").Value = "'  -0.49%  "

# Row 16
$ws.Range("D16").Value = "'29.252.45"
$ws.Range("E16").Value = "'  -0.29%  "

# Row 17
$ws.Range("E17").Value = "'  +1.05%  "

# Row 18
$ws.Range("D18").Value = "'243.39"
$ws.Range("E18").Value = "'  -0.71%  "

# Row 19
$ws.Range("D19").Value = "'2.141.14"
$ws.Range("E19").Value = "'  +2.18%  "

# Row 20
$ws.Range("D20").Value = "'0.000007812"
$ws.Range("E20").Value = "'  -1.08%  "

# Row 21
$ws.Range("E21").Value = "'  -1.11%  "

# Row 22
$ws.Range("D22").Value = "'0.9996"
$ws.Range("E22").Value = "'  +0.08%  "

# Row 23
$ws.Range("D23").Value = "'7.922"
$ws.Range("E23").Value = "'  -0.37%  "

# Row 24
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "'  +0.07%  "

# Row 25
$ws.Range("D25").Value = "'0.1613"
$ws.Range("E25").Value = "'  +1.72%  "

# Row 26
$ws.Range("D26").Value = "'162.75"
$ws.Range("E26").Value = "'  -0.74%  "

# Row 27
$ws.Range("D27").Value = "'8.906"
$ws.Range("E27").Value = "'  -1.41%  "

# Row 28
$ws.Range("D28").Value = "'18.59"
$ws.Range("E28").Value = "'  +1.32%  "

# Row 29
$ws.Range("D29").Value = "'1.353"
$ws.Range("E29").Value = "'  -0.82%  "

# Row 30
$ws.Range("D30").Value = "'1.500"
$ws.Range("E30").Value = "'  +0.97%  "

# Row 31
$ws.Range("D31").Value = "'4.428"
$ws.Range("E31").Value = "'  +0.86%  "

# Row 32
$ws.Range("D32").Value = "'4.264"
$ws.Range("E32").Value = "'  +2.68%  "

# Row 33
$ws.Range("D33").Value = "'0.05177"
$ws.Range("E33").Value = "'  -1.96%  "

# Row 34
$ws.Range("D34").Value = "'0.8154"
$ws.Range("E34").Value = "'  +12.92%  "

# Row 35
$ws.Range("D35").Value = "'1.931"
$ws.Range("E35").Value = "'  -0.58%  "

# Row 36
$ws.Range("D36").Value = "'1.174"
$ws.Range("E36").Value = "'  -2.26%  "

# Row 37
$ws.Range("D37").Value = "'2.684"
$ws.Range("E37").Value = "'  +0.26%  "

# Row 38
$ws.Range("E38").Value = "'  -0.46%  "

# Row 39
$ws.Range("E39").Value = "'  -1.04%  "

# Row 40
$ws.Range("D40").Value = "'1.159.71"
$ws.Range("E40").Value = "'  -6.58%  "

# Row 41
$ws.Range("D41").Value = "'6.208"
$ws.Range("E41").Value = "'  +2.00%  "

# Row 42
$ws.Range("D42").Value = "'0.8943"
$ws.Range("E42").Value = "'  -1.52%  "

# Row 43
$ws.Range("D43").Value = "'72.79"
$ws.Range("E43").Value = "'  -0.47%  "

# Row 44
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "'  +0.01%  "

# Row 45
$ws.Range("D45").Value = "'101.81"
$ws.Range("E45").Value = "'  -1.72%  "

# Row 46
$ws.Range("D46").Value = "'2.035.61"
$ws.Range("E46").Value = "'  +1.41%  "

# Row 47
$ws.Range("D47").Value = "'0.5181"
$ws.Range("E47").Value = "'  -2.77%  "

# Row 48
$ws.Range("D48").Value = "'1.783"
$ws.Range("E48").Value = "'  +1.00%  "

# Row 49
$ws.Range("D49").Value = "'9.361"
$ws.Range("E49").Value = "'  +0.93%  "

# Row 50
$ws.Range("B50").Value = "'Aptos"
$ws.Range("C50").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'7.067"
$ws.Range("E50").Value = "'  -0.21%  "

# Row 51
$ws.Range("B51").Value = "'TheSandbox"
$ws.Range("C51").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "'0.4280"
$ws.Range("E51").Value = "'  -1.25%  "

